# Add the new Twitter-bot error message as a new row (A3) below the
# existing two rows, matching the author's new "already favorited" error.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The text begins with a literal apostrophe. Assigning it straight to
# .Value/.Value2 makes Excel treat the leading "'" as a text-quote-prefix
# indicator and stamp the cell with a new (quotePrefix) style, which the
# original workbook never had. Routing the literal through a text formula
# first, then converting the formula to a plain value via copy/paste-values,
# stores it as an ordinary shared-string cell with no style change at all.
$ws.Range("A3").Formula = "=""'message': 'You have already favorited this status.', 'code': 139}"""
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("A3").Select()
